# The original workbook stores every cell (even ones that look numeric,
# like "0.207" or "-0.015") as TEXT (shared-string) values rather than as
# real numbers - this is a balance-table export where every "number" is
# really a formatted string. To update the five corrected figures while
# keeping them as text cells (not reinterpreted as numeric values by
# Excel's normal type-inference), we temporarily mark each target cell as
# Text before typing the new value, then strip the temporary number
# format back off again (via PasteSpecial of formats copied from an
# untouched, default-styled cell) so the cells end up back on the
# workbook's default/general style, matching the rest of the sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$targets = @{
    "E4"  = "0.193"     # prop_independiente_trabajando, group 1 mean: 0.192 -> 0.193
    "G4"  = "-0.014"    # prop_independiente_trabajando, mean difference: -0.015 -> -0.014
    "E6"  = "0.200"     # prop_independiente_total, group 1 mean: 0.199 -> 0.200
    "E10" = "0.673"     # prop_formal_no_indep, group 1 mean: 0.675 -> 0.673
    "G10" = "0.163***"  # prop_formal_no_indep, mean difference: 0.164*** -> 0.163***
}

foreach ($addr in $targets.Keys) {
    $ws.Range($addr).NumberFormat = "@"
    $ws.Range($addr).Value = $targets[$addr]
}

# Restore the plain/default formatting on the edited cells so they match
# the look of the rest of the (untouched) table instead of staying
# flagged with an explicit "Text" number format.
$ws.Range("A1").Copy()
foreach ($addr in $targets.Keys) {
    $ws.Range($addr).PasteSpecial(-4122)
}
$excel.CutCopyMode = 0
